$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.204.71"
$ws.Range("E2").Value = "  -2.93%  "

$ws.Range("D3").Value = "3.299.60"
$ws.Range("E3").Value = "  -3.42%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "'557.85"
$ws.Range("E5").Value = "  -3.19%  "

$ws.Range("D6").Value = "'141.65"
$ws.Range("E6").Value = "  -4.67%  "

$ws.Range("D8").Value = "3.298.89"
$ws.Range("E8").Value = "  -3.48%  "

$ws.Range("D9").Value = "'0.472"
$ws.Range("E9").Value = "  -2.40%  "

$ws.Range("D10").Value = "'7.84"
$ws.Range("E10").Value = "  -1.65%  "

$ws.Range("E11").Value = "  -4.04%  "

$ws.Range("E12").Value = "  -2.17%  "

$ws.Range("D13").Value = "3.873.08"
$ws.Range("E13").Value = "  -3.20%  "

$ws.Range("E14").Value = "  +0.45%  "

$ws.Range("E15").Value = "  -5.45%  "

$ws.Range("D16").Value = "3.291.89"
$ws.Range("E16").Value = "  -3.44%  "

$ws.Range("E17").Value = "  -3.39%  "

$ws.Range("D18").Value = "60.219.38"
$ws.Range("E18").Value = "  -2.86%  "

$ws.Range("E19").Value = "  -3.44%  "

$ws.Range("D20").Value = "'14.39"
$ws.Range("E20").Value = "  -0.88%  "

$ws.Range("D21").Value = "'8.61"
$ws.Range("E21").Value = "  -4.18%  "

$ws.Range("D22").Value = "'373.71"
$ws.Range("E22").Value = "  -1.95%  "

$ws.Range("D23").Value = "'74.22"
$ws.Range("E23").Value = "  -0.92%  "

$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  +0.00%  "

$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").Value = "'0.541"
$ws.Range("E25").Value = "  -4.86%  "

$ws.Range("D26").Value = "3.445.76"
$ws.Range("E26").Value = "  -3.33%  "

$ws.Range("E27").Value = "  -7.95%  "

$ws.Range("E28").Value = "  -4.42%  "

$ws.Range("D29").Value = "'0.997"
$ws.Range("E29").Value = "  -0.19%  "

$ws.Range("D30").Value = "'7.23"
$ws.Range("E30").Value = "  -5.15%  "

$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("D32").Value = "'7.60"
$ws.Range("E32").Value = "  -3.97%  "

$ws.Range("E33").Value = "  -4.00%  "

$ws.Range("D34").Value = "'22.52"
$ws.Range("E34").Value = "  -2.36%  "

$ws.Range("D35").Value = "'1.26"
$ws.Range("E35").Value = "  -5.78%  "

$ws.Range("E36").Value = "  -6.27%  "

$ws.Range("D37").Value = "'165.77"
$ws.Range("E37").Value = "  -2.36%  "

$ws.Range("D38").Value = "'1.52"
$ws.Range("E38").Value = "  -5.70%  "

$ws.Range("D39").Value = "'6.71"
$ws.Range("E39").Value = "  -2.81%  "

$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "'26.86"
$ws.Range("E40").Value = "  -11.95%  "

$ws.Range("B41").Value = "RenzoRestakedETH"
$ws.Range("C41").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D41").Value = "3.334.84"
$ws.Range("E41").Value = "  -3.35%  "

$ws.Range("D42").Value = "'0.0737"
$ws.Range("E42").Value = "  -5.59%  "

$ws.Range("D43").Value = "'41.93"
$ws.Range("E43").Value = "  -1.26%  "

$ws.Range("D44").Value = "'0.751"
$ws.Range("E44").Value = "  -3.29%  "

$ws.Range("D45").Value = "'4.18"
$ws.Range("E45").Value = "  -4.38%  "

$ws.Range("E46").Value = "  -5.37%  "

$ws.Range("E47").Value = "  -4.81%  "

$ws.Range("D48").Value = "2.368.85"
$ws.Range("E48").Value = "  -6.86%  "

$ws.Range("E49").Value = "  +0.18%  "

$ws.Range("E50").Value = "  -6.12%  "

$ws.Range("E51").Value = "  -6.73%  "
